$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.553.37'
$ws.Range('E2').Value = '  +2.34%  '
$ws.Range('D3').Value = '1.670.06'
$ws.Range('E3').Value = '  +1.68%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9986'
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '238.83'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9998'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.18%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4812'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('E8').Value = '  +1.68%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06190'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.18%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07011'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.62%  '
$ws.Range('D11').Value = '1.668.68'
$ws.Range('E11').Value = '  +1.58%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.90'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.41%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.5911'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.386'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.08%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '75.22'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.31%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9996'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.05%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.9992'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.09%  '
$ws.Range('D18').Value = '25.540.09'
$ws.Range('E18').Value = '  +2.25%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000006779'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.56%  '
$ws.Range('E20').Value = '  +0.83%  '
$ws.Range('D21').Value = '1.882.00'
$ws.Range('E21').Value = '  +1.41%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.450'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.743'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.78%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.288'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '136.75'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.63%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.05'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.71%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.392'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.15%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.727'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.07%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '104.86'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.98%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.976'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +6.40%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.07811'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.04%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.660'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.9986'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04247'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.04%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.614'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.12%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6101'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.21%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9522'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.56%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.591'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.79%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.8586'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.41%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9994'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.08%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.859'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.89%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.01483'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.48%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '95.93'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.40%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3779'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.70%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.844'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.59%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1120'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.60%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.214'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.92%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.05250'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.07%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '29.83'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.37%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.356'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.82%  '
$ws.Range('B51').Value = 'TrueUSD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.001'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.19%  '
